$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Efnb2", "Ephb4", "ECs", 3, 1, 37.98277566666666, 113.948327, 0.697850645410475, 0.6978506454104751, 3, 1, 49.146613, 147.439839, 0.8214902885327882, 0.8214902885327882, 1866.724776355484, 16800.52298719935, 0.5732775280510436, 0.5732775280510437),
  @("ECs", "Efnb2", "Ephb4", "FAPs", 3, 1, 37.98277566666666, 113.948327, 0.697850645410475, 0.6978506454104751, 3, 1, 8.144916333333333, 24.434749, 0.1361430475126621, 0.1361430475126621, 309.3665299127692, 2784.298769214923, 0.09500751357486023, 0.09500751357486024),
  @("ECs", "Efnb2", "Ephb4", "MuSCs", 3, 1, 37.98277566666666, 113.948327, 0.697850645410475, 0.6978506454104751, 3, 1, 2.534635, 7.603904999999999, 0.04236666395454969, 0.04236666395454969, 96.27247260188165, 866.4522534169348, 0.02956560378457121, 0.02956560378457121),
  @("FAPs", "Efnb2", "Ephb4", "ECs", 3, 1, 9.685730333333334, 29.057191, 0.1779541659542351, 0.1779541659542352, 3, 1, 49.146613, 147.439839, 0.8214902885327882, 0.8214902885327882, 476.0208403146944, 4284.18756283225, 0.1461876191353563, 0.1461876191353563),
  @("FAPs", "Efnb2", "Ephb4", "FAPs", 3, 1, 9.685730333333334, 29.057191, 0.1779541659542351, 0.1779541659542352, 3, 1, 8.144916333333333, 24.434749, 0.1361430475126621, 0.1361430475126621, 78.88946319222877, 710.005168730059, 0.02422722247058359, 0.0242272224705836),
  @("FAPs", "Efnb2", "Ephb4", "MuSCs", 3, 1, 9.685730333333334, 29.057191, 0.1779541659542351, 0.1779541659542352, 3, 1, 2.534635, 7.603904999999999, 0.04236666395454969, 0.04236666395454969, 24.54979110342833, 220.948119930855, 0.007539324348295247, 0.007539324348295248),
  @("Inflammatory-Mac", "Efnb2", "Ephb4", "ECs", 3, 1, 0.5676613333333332, 1.702984, 0.01042953867610283, 0.01042953867610283, 3, 1, 49.146613, 147.439839, 0.8214902885327882, 0.8214902885327882, 27.89863186439733, 251.087686779576, 0.008567764736295587, 0.008567764736295589),
  @("Inflammatory-Mac", "Efnb2", "Ephb4", "FAPs", 3, 1, 0.5676613333333332, 1.702984, 0.01042953867610283, 0.01042953867610283, 3, 1, 8.144916333333333, 24.434749, 0.1361430475126621, 0.1361430475126621, 4.623554065668444, 41.611986591016, 0.001419909179515815, 0.001419909179515815),
  @("Inflammatory-Mac", "Efnb2", "Ephb4", "MuSCs", 3, 1, 0.5676613333333332, 1.702984, 0.01042953867610283, 0.01042953867610283, 3, 1, 2.534635, 7.603904999999999, 0.04236666395454969, 0.04236666395454969, 1.438814283613333, 12.94932855252, 0.0004418647602914277, 0.0004418647602914277),
  @("MuSCs", "Efnb2", "Ephb4", "ECs", 3, 1, 5.823095333333334, 17.469286, 0.1069866739681064, 0.1069866739681064, 3, 1, 49.146613, 147.439839, 0.8214902885327882, 0.8214902885327882, 286.1854128094394, 2575.668715284954, 0.08788851366722308, 0.08788851366722308),
  @("MuSCs", "Efnb2", "Ephb4", "FAPs", 3, 1, 5.823095333333334, 17.469286, 0.1069866739681064, 0.1069866739681064, 3, 1, 8.144916333333333, 24.434749, 0.1361430475126621, 0.1361430475126621, 47.42862429102378, 426.857618619214, 0.01456549183726161, 0.01456549183726161),
  @("MuSCs", "Efnb2", "Ephb4", "MuSCs", 3, 1, 5.823095333333334, 17.469286, 0.1069866739681064, 0.1069866739681064, 3, 1, 2.534635, 7.603904999999999, 0.04236666395454969, 0.04236666395454969, 14.75942124020333, 132.83479116183, 0.004532668463621734, 0.004532668463621734),
  @("Resolving-Mac", "Efnb2", "Ephb4", "ECs", 3, 1, 0.3689676666666666, 1.106903, 0.006778975991080511, 0.006778975991080512, 3, 1, 49.146613, 147.439839, 0.8214902885327882, 0.8214902885327882, 18.13351112317967, 163.201600108617, 0.005568862942869572, 0.005568862942869573),
  @("Resolving-Mac", "Efnb2", "Ephb4", "FAPs", 3, 1, 0.3689676666666666, 1.106903, 0.006778975991080511, 0.006778975991080512, 3, 1, 8.144916333333333, 24.434749, 0.1361430475126621, 0.1361430475126621, 3.005210774705222, 27.046896972347, 0.0009229104504408698, 0.0009229104504408699),
  @("Resolving-Mac", "Efnb2", "Ephb4", "MuSCs", 3, 1, 0.3689676666666666, 1.106903, 0.006778975991080511, 0.006778975991080512, 3, 1, 2.534635, 7.603904999999999, 0.04236666395454969, 0.04236666395454969, 0.9351983618016665, 8.416785256214999, 0.0002872025977700684, 0.0002872025977700684)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $rowIndex = $i + 2
  for ($j = 0; $j -lt $data[$i].Count; $j++) {
    $colIndex = $j + 1
    $ws.Cells.Item($rowIndex, $colIndex).Value = $data[$i][$j]
  }
}

$ws.Rows("17:21").Delete()